$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns from ORG_FAC_* to ORG_FACT_*
$ws.Range("A1").Value = "ORG_FACT_PK"
$ws.Range("B1").Value = "ORG_FACT_NAME"
$ws.Range("C1").Value = "ORG_FACT_ABBR_NAME"
$ws.Range("D1").Value = "ORG_FACT_IDENOLD"
$ws.Range("E1").Value = "ORG_FACT_IDENNEW"
$ws.Range("F1").Value = "ORG_FACT_STATUS"

# Update the active selection to F1
$ws.Range("F1").Select()
